$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Subject completes response to perturbation having steered the vehicle back to the center of the lane. Normally this would be tagged with temporal scope but avoiding definitions here."
$ws.Range("D5").Select()
